$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Odd_Over25_FT / Odd_Under25_FT
$ws.Range("Q2").Value = 1.67
$ws.Range("R2").Value = 2.15

# Row 3: Odd_Over25_FT / Odd_Under25_FT
$ws.Range("Q3").Value = 2.5
$ws.Range("R3").Value = 1.5
